$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'3.68%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-7.78%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.229"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.39%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'1.94%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.715"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.76%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'MXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'0.8637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.67%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'FTXToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'0.9635"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'12.14%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'One"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.01056"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1,666.80%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1411"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.86%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.14%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03179"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.95%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.43%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001549"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.65%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.005810"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.73%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'LEO"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.501"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.34%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.218"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.02%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.221"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.94%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3178"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03477"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.51%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-0.67%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.533"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.54%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'1.20%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-2.11%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.18%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004800"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'15.42%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'1.19%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03814"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.64%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'KickToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.005670"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.55%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.1100"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.89%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'9.54%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01072"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.29%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005238"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.51%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'40.86%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002130"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-13.86%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").Style = "Normal"
